# Rename the workbook's only worksheet from "Sheet1" to "packages",
# matching the change reflected in xl/workbook.xml (<sheet name="packages" .../>).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "packages"
